# Update forests data - 2025-10-08 12:18
#
# The 4 listings currently on the "New" sheet (rows 2-5) are the
# previously-scraped items that now graduate onto the "Previously added"
# sheet (appended as rows 161-164, with their hyperlinks). The "New" sheet
# is then repopulated with 3 freshly scraped listings (rows 2-4).

$wb     = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# A never-touched data row whose formatting (s="3"/"4"/"2" per column) we
# reuse as a style template everywhere below - both sheets share the same
# styles.xml, so copying across sheets is safe.
$styleTemplate = $wsPrev.Range("A160:F160")

# ---------------------------------------------------------------------------
# Data that is moving from "New" -> "Previously added" (was rows 2-5 on
# "New"), in order: link, price, district, area, cadastre, date.
# ---------------------------------------------------------------------------
$movedRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/kuku-pag/mfhkk.html",
      "45 000 €", "Jēkabpils un raj.", "6 ha.", "56700050044", 45936.663888888885),
    @("https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kaplavas-pag/dmief.html",
      "126 000 €", "Krāslava un raj.", "20 ha.", "60700090144", 45937.475694444445),
    @("https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/purenu-pag/ihhdc.html",
      "70 000 €", "Ludza un raj.", "7 ha.", "68880010168", 45937.37430555555),
    @("https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/gibulu-pag/demgx.html",
      "19 000 €", "Talsi un raj.", "7 ha.", "", 45937.46666666667)
)

# ---------------------------------------------------------------------------
# Freshly scraped listings that now populate "New" rows 2-4.
# ---------------------------------------------------------------------------
$newRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/kubulu-pag/cddkb.html",
      "50 000 €", "Balvi un raj.", "6 ha.", "38580030185", 45938.54236111111),
    @("https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/malinovas-pag/mknjo.html",
      "20 500 €", "Daugavpils un raj.", "3 ha.", "44700010164", 45937.643055555556),
    @("https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/rusonas-pag/hfiff.html",
      "15 000 €", "Preiļi un raj.", "1 ha.", "7670 002 0088", 45938.44861111111)
)

function Set-TextCell($range, [string]$text) {
    # Prefix with an apostrophe so numeric/blank-looking strings (e.g. a
    # cadastre number) are stored as text rather than being auto-converted
    # to a number by the COM layer's "smart typing" behaviour.
    $range.Value = "'" + $text
}

function Set-RowData($ws, [int]$r, $row) {
    Set-TextCell $ws.Range("A$r") $row[0]
    Set-TextCell $ws.Range("B$r") $row[1]
    Set-TextCell $ws.Range("C$r") $row[2]
    Set-TextCell $ws.Range("D$r") $row[3]
    Set-TextCell $ws.Range("E$r") $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# ---------------------------------------------------------------------------
# Append the moved rows to the bottom of "Previously added" (rows 161-164).
# ---------------------------------------------------------------------------
$lastRow = $wsPrev.UsedRange.Rows.Count

for ($i = 0; $i -lt $movedRows.Count; $i++) {
    $destRow = $lastRow + 1 + $i
    $dstRange = $wsPrev.Range("A${destRow}:F$destRow")

    # 1) Copy formatting/styles down onto the new row.
    $styleTemplate.Copy($dstRange)

    # 2) Register the hyperlink (this also applies Excel's built-in
    #    "Hyperlink" style to the cell).
    $row = $movedRows[$i]
    $wsPrev.Hyperlinks.Add($wsPrev.Range("A$destRow"), $row[0])

    # 3) Re-copy the formatting to restore this sheet's own link style; the
    #    hyperlink relationship added in step 2 is unaffected.
    $styleTemplate.Copy($dstRange)

    # 4) Finally, write the real values into place.
    Set-RowData $wsPrev $destRow $row
}

# ---------------------------------------------------------------------------
# Wipe the old hyperlinks on "New" (removes all hyperlink relationships on
# the sheet in one go) and delete the now-vacated row 5.
# ---------------------------------------------------------------------------
$wsNew.Range("A2").Hyperlinks.Delete()
$wsNew.Rows.Item(5).Delete()

# ---------------------------------------------------------------------------
# Populate "New" rows 2-4 with the freshly scraped listings.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]
    $dstRange = $wsNew.Range("A${r}:F$r")

    $wsNew.Hyperlinks.Add($wsNew.Range("A$r"), $row[0])
    # Restore normal (non-hyperlink) styling after Add() applied the
    # built-in "Hyperlink" style to column A.
    $styleTemplate.Copy($dstRange)

    Set-RowData $wsNew $r $row
}

Write-Output "Edit complete"
